$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their text (string) storage, matching the
# original inline-string cell type, so numeric-looking values like "1.009"
# are not auto-converted into numbers by Excel.
$ws.Columns("D:E").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '20.401.47'
$ws.Cells.Item(2, 5).Value = '  -7.15%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.442.84'
$ws.Cells.Item(3, 5).Value = '  -7.09%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.009'
$ws.Cells.Item(4, 5).Value = '  +0.53%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '1.008'
$ws.Cells.Item(5, 5).Value = '  +0.66%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '277.46'
$ws.Cells.Item(6, 5).Value = '  -4.45%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.3717'
$ws.Cells.Item(7, 5).Value = '  -5.26%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.3092'
$ws.Cells.Item(8, 5).Value = '  -4.17%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '41.01'
$ws.Cells.Item(9, 5).Value = '  -7.76%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '1.007'
$ws.Cells.Item(10, 5).Value = '  -5.98%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '0.06547'
$ws.Cells.Item(11, 5).Value = '  -8.97%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.52%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '5.415'
$ws.Cells.Item(13, 5).Value = '  -4.29%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '17.19'
$ws.Cells.Item(14, 5).Value = '  -7.63%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'WrappedEther'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(15, 4).Value = '1.452.86'
$ws.Cells.Item(15, 5).Value = '  -6.65%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).Value = '6.131'
$ws.Cells.Item(16, 5).Value = '  -7.59%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '0.00001014'
$ws.Cells.Item(17, 5).Value = '  -8.49%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '76.93'
$ws.Cells.Item(18, 5).Value = '  -7.67%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '0.05816'
$ws.Cells.Item(19, 5).Value = '  -11.85%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.64%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '5.718'
$ws.Cells.Item(21, 5).Value = '  -8.07%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '14.45'
$ws.Cells.Item(22, 5).Value = '  -6.27%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '10.97'
$ws.Cells.Item(23, 5).Value = '  -2.32%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '2.299'
$ws.Cells.Item(24, 5).Value = '  -2.62%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '20.408.54'
$ws.Cells.Item(25, 5).Value = '  -7.21%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26, 4).Value = '142.93'
$ws.Cells.Item(26, 5).Value = '  -3.15%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(27, 4).Value = '2.218'
$ws.Cells.Item(27, 5).Value = '  -7.08%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '17.01'
$ws.Cells.Item(28, 5).Value = '  -8.40%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '1.617.71'
$ws.Cells.Item(29, 5).Value = '  -6.85%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '108.85'
$ws.Cells.Item(30, 5).Value = '  -8.65%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '0.9178'
$ws.Cells.Item(31, 5).Value = '  -6.23%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '5.381'
$ws.Cells.Item(32, 5).Value = '  -8.68%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '3.386'
$ws.Cells.Item(33, 5).Value = '  -30.37%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '0.07735'
$ws.Cells.Item(34, 5).Value = '  -6.77%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '8.332'
$ws.Cells.Item(35, 5).Value = '  -8.52%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Frax'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(36, 4).Value = '1.008'
$ws.Cells.Item(36, 5).Value = '  +0.68%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Aptos'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(37, 4).Value = '10.92'
$ws.Cells.Item(37, 5).Value = '  +2.22%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '1.411'
$ws.Cells.Item(38, 5).Value = '  -12.43%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.05611'
$ws.Cells.Item(39, 5).Value = '  -6.43%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '4.729'
$ws.Cells.Item(40, 5).Value = '  -7.30%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '1.125'
$ws.Cells.Item(41, 5).Value = '  -6.51%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.02037'
$ws.Cells.Item(42, 5).Value = '  -9.58%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '0.1908'
$ws.Cells.Item(43, 5).Value = '  -7.03%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '3.591'
$ws.Cells.Item(44, 5).Value = '  -4.23%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.5317'
$ws.Cells.Item(45, 5).Value = '  -8.05%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '12.18'
$ws.Cells.Item(46, 5).Value = '  -6.34%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '0.5145'
$ws.Cells.Item(47, 5).Value = '  -7.26%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '110.79'
$ws.Cells.Item(48, 5).Value = '  -5.49%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '1.762'
$ws.Cells.Item(49, 5).Value = '  -5.96%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '1.055'
$ws.Cells.Item(50, 5).Value = '  -7.04%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '1.007'
$ws.Cells.Item(51, 5).Value = '  +0.39%  '
